# Update the "Förändrad" (Changed) date column (C) for rows 2 through 28
# from 2024-06-16 (serial 45459) to 2024-06-17 (serial 45460), i.e. +1 day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    $cell.Value = $current.AddDays(1)
}
